$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.531.94"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.952.13"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.49"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.23"
$ws.Range("E7").Value = "  +5.18%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +4.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0790"
$ws.Range("E10").Value = "  -6.89%  "
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.23"
$ws.Range("E12").Value = "  +6.65%  "
$ws.Range("D13").Value = "2.239.64"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.826"
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.48"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.25"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "1.952.26"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "36.424.57"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.30"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "0.0₃0852"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.08"
$ws.Range("E21").Value = "  +2.66%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.93"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("E25").Value = "  +2.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.144"
$ws.Range("E26").Value = "  +7.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.16"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.10"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.30"
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.32"
$ws.Range("E30").Value = "  +21.44%  "
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.74"
$ws.Range("E32").Value = "  +4.06%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.46"
$ws.Range("E34").Value = "  +7.45%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.28"
$ws.Range("E35").Value = "  +5.22%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.41"
$ws.Range("E37").Value = "  +9.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("E39").Value = "  -11.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0968"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0210"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.76"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "1.361.16"
$ws.Range("E45").Value = "  +2.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.60"
$ws.Range("E46").Value = "  +3.66%  "
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.12"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.19"
$ws.Range("E50").Value = "  +5.02%  "
$ws.Range("D51").Value = "2.135.31"
$ws.Range("E51").Value = "  +1.08%  "
